$d = $word.ActiveDocument

# 1. Merge the spell-checked "flutter_blue" run back into plain text (removes proofErr spellcheck markers)
$d.Content.Find.Execute("called flutter_blue that", $true, $false, $false, $false, $false, $true, 1, $false, "called flutter_blue that", 2) | Out-Null

# 2. Append a placeholder empty paragraph at the end of the document, then replace its
#    range with the full Week 16/17/18 OOXML fragment (keeps "On week 15..." paragraph intact).
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.InsertParagraphAfter()
$placeholder = $d.Paragraphs($d.Paragraphs.Count)

$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Week 16</w:t></w:r></w:p><w:p><w:r><w:t>This week was a week where a lot of time was put into creating the structure for the Bluetooth. After looking at the flutter blue repository I found many different files that we needed in our model part of the app. Since we are working in a model-view-</w:t></w:r><w:r><w:t>view Model</w:t></w:r><w:r><w:t xml:space="preserve"> we need to create the dependencies in separate folders and then initialize them in the </w:t></w:r><w:r><w:t>View Model and control it before using them in the view</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>We are still some ways from being done with the mower but hopefully, it goes fast. Joakim finished the first part of the Joystick and will now begin looking at the view for the Bluetooth list.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Week 17</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The goal for this week has been to hopefully finish the Bluetooth part. I and Joakim have been working together with creating the view model and view to make sure that everything works. Joakim finished all the separate parts for the view and my task is now to </w:t></w:r><w:r><w:t xml:space="preserve">make it work with the view model. At the end of the week, we were ready to try the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Bluetooth</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> but the Bluetooth sadly didn’t work. It seems that we have a problem with the permissions for the project. Since we use Bluetooth, we need to ask for permission before it is used but there seems to be something wrong.</w:t></w:r><w:r><w:t xml:space="preserve"> Will continue to work on it next week.</w:t></w:r><w:r><w:t xml:space="preserve"> Since Joakim is done </w:t></w:r><w:r><w:t xml:space="preserve">with </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>view</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> we decided that he will start on the </w:t></w:r><w:r><w:t>HTTP</w:t></w:r><w:r><w:t xml:space="preserve"> requests from the REST API</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Week 18</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">For this week’s mission, it was to make sure that the Bluetooth and the rest API requests were working. We in the group had set a target that we should be done with all the different connections this week to make sure that we are at pace with the project. The problem that we still had with Bluetooth was that the permissions were wrong. We sat with this problem for a while and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>looked into</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> other libraries and examples of code to see what we could have done wrong. The only output we received was that there were missing permissions in the manifest. This means that the android project isn’t set up right. At the end of the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>week</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> we realized that where we are looking for permissions, we had forgotten to ask for all the required permissions. The earlier belief was that we only needed to ask for Bluetooth when it in fact was needed to check three different permissions. With this done we could finally find different devices on our app. The new problem was that we couldn’t find the mower Bluetooth device. </w:t></w:r><w:r><w:t xml:space="preserve">After some research and testing with integrating a search filter for the mower, we realized that there might be a problem with that the library that is being used is using BLE instead of regular Bluetooth. I will contact Andreas about receiving a raspberry-pi 4 instead and see if that solves the issue and hopefully, we can be done with Bluetooth next week. </w:t></w:r></w:p>
'@

$placeholder.Range.InsertXML($newContentXml) | Out-Null

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
